$wb = $excel.ActiveWorkbook

# --- Update the "Share of Cost Effective Capacity Built in a Single Year"
# sheet: row 7 ("onshore wind es") changes from 0.33 to 0.2 across all
# year columns (B:AE) ---
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$ws.Range("B7:AE7").Value = 0.2

# --- Make this sheet the active/selected tab, with the selection resting
# on the row that was just edited ---
$ws.Activate()
$ws.Range("B7:AE7").Select()
